# Sync updated BOM for multitarget
#
# - Move "Digikey Final" to be the first (leftmost) tab.
# - Rename "Digikey Upload" to "DO NOT USE" (formulas on other sheets that
#   reference it by name are retargeted automatically).
# - Make "From Eagle" the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Reorder: put "Digikey Final" before whatever is currently the first sheet.
$wb.Worksheets.Item("Digikey Final").Move($wb.Worksheets.Item(1))

# Flag the old "Digikey Upload" sheet as deprecated.
$wb.Worksheets.Item("Digikey Upload").Name = "DO NOT USE"

# Select "From Eagle" as the active sheet/tab.
$wb.Worksheets.Item("From Eagle").Activate()
